$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Skill")
$ws2 = $wb.Worksheets.Item("디스크립션")

# Insert a new column before column F (abnormal), shifting abnormal -> G, type -> H
$ws1.Range("F1").EntireColumn.Insert()

# New header for the inserted column F
$ws1.Range("F1").Value = "cooltime"

# Fill cooltime values for rows 2-7
$ws1.Range("F2:F7").Value = 0.1

# Column width for the new F (cooltime) column = 24.25 chars, no bestFit.
# (Column G, the shifted former "abnormal" column, keeps its original 15.5/
# bestFit width automatically via the Insert() shift above - leave it alone.)
$ws1.Range("F1").EntireColumn.ColumnWidth = 23.5

# Selection changes per diff (set sheet2's selection first, then re-activate
# sheet1 last so it remains the tab-selected sheet, matching the diff)
$ws2.Range("D29").Select() | Out-Null
$ws1.Range("E16").Select() | Out-Null
